$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "clean" numeric-looking string that would
# otherwise be auto-converted to a Number (losing e.g. trailing zeros) -
# pre-format as Text so they round-trip as the exact literal string.
$forceTextCells = @("D4", "D5", "D22", "D26", "D40", "D41")
foreach ($c in $forceTextCells) {
  $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.406.94'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '1.869.38'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '243.50'
$ws.Range("E6").Value = '  -2.34%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '0.07946'
$ws.Range("E8").Value = '  -0.57%  '
$ws.Range("D9").Value = '0.3134'
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("D10").Value = '24.48'
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("D11").Value = '0.07812'
$ws.Range("E11").Value = '  -4.92%  '
$ws.Range("D12").Value = '1.898.34'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("D13").Value = '93.74'
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").Value = '5.172'
$ws.Range("D15").Value = '0.7023'
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").Value = '6.502'
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").Value = '0.000008521'
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").Value = '29.463.00'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("E19").Value = '  +3.56%  '
$ws.Range("D20").Value = '2.144.72'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '7.659'
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").Value = '0.1549'
$ws.Range("E25").Value = '  -2.99%  '
$ws.Range("D26").Value = '9.010'
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").Value = '161.66'
$ws.Range("E28").Value = '  +1.58%  '
$ws.Range("D29").Value = '1.507'
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").Value = '4.312'
$ws.Range("E30").Value = '  -2.09%  '
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("D32").Value = '1.216'
$ws.Range("E32").Value = '  +2.14%  '
$ws.Range("D33").Value = '0.05267'
$ws.Range("D34").Value = '1.902'
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").Value = '0.7574'
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D37").Value = '2.708'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.284.89'
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01877'
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").Value = '2.770'
$ws.Range("E40").Value = '  +0.67%  '
$ws.Range("D41").Value = '0.8960'
$ws.Range("E41").Value = '  -1.64%  '
$ws.Range("D42").Value = '109.58'
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D43").Value = '6.017'
$ws.Range("E43").Value = '  -6.61%  '
$ws.Range("D44").Value = '70.97'
$ws.Range("E44").Value = '  -4.23%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000127'
$ws.Range("E46").Value = '  -3.22%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '2.040.11'
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("D48").Value = '1.806'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D49").Value = '9.603'
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("D50").Value = '0.5187'
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("D51").Value = '0.4297'
$ws.Range("E51").Value = '  -0.93%  '
